$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()

$ws.Range("H137").Value = 2441324.8
$ws.Range("I137").Value = 5264921.5
$ws.Range("J137").Value = 2763.6365
$ws.Range("K137").Value = 15794764.5
$ws.Range("L137").Value = 8290.9095
$ws.Range("M137").Value = -15792214.5
$ws.Range("N137").Value = -13390.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4291578.5
$ws.Range("I32").Value = 4778590
$ws.Range("K32").Value = 4778590
$ws.Range("M32").Value = -4778303

$ws.Range("H61").Value = 40080804
$ws.Range("I61").Value = 45500844
$ws.Range("J61").Value = 333865
$ws.Range("K61").Value = 45500844
$ws.Range("L61").Value = 333865
$ws.Range("M61").Value = -45500632
$ws.Range("N61").Value = -334289

$ws.Range("H74").Value = 7637478
$ws.Range("I74").Value = 11953229
$ws.Range("J74").Value = 84913.75
$ws.Range("K74").Value = 11953229
$ws.Range("L74").Value = 84913.75
$ws.Range("M74").Value = -11952355
$ws.Range("N74").Value = -86661.75

$ws.Range("H77").Value = 7637478
$ws.Range("I77").Value = 11953229
$ws.Range("J77").Value = 84913.75
$ws.Range("K77").Value = 59766145
$ws.Range("L77").Value = 424568.75
$ws.Range("M77").Value = -59761777
$ws.Range("N77").Value = -433304.75

$ws.Range("H132").Value = 74167.82
$ws.Range("I132").Value = 79457.234
$ws.Range("J132").Value = 69583.664
$ws.Range("K132").Value = 238371.702
$ws.Range("L132").Value = 208750.992
$ws.Range("M132").Value = -235841.702
$ws.Range("N132").Value = -213810.992

$ws.Range("H136").Value = 40080804
$ws.Range("I136").Value = 45500844
$ws.Range("J136").Value = 333865
$ws.Range("K136").Value = 136502532
$ws.Range("L136").Value = 1001595
$ws.Range("M136").Value = -136499982
$ws.Range("N136").Value = -1006695

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7944.476
$ws.Range("I31").Value = 30348.475
$ws.Range("J31").Value = 1395.6154
$ws.Range("K31").Value = 30348.475
$ws.Range("L31").Value = 1395.6154
$ws.Range("M31").Value = -30053.475
$ws.Range("N31").Value = -1985.6154

$ws.Range("H34").Value = 7944.476
$ws.Range("I34").Value = 30348.475
$ws.Range("J34").Value = 1395.6154
$ws.Range("K34").Value = 30348.475
$ws.Range("L34").Value = 1395.6154
$ws.Range("M34").Value = -30146.475
$ws.Range("N34").Value = -1799.6154

$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -27246

$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -86232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 263.33334
$ws.Range("I10").Value = 263.33334
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 790.0000200000001
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -651.0000200000001
$ws.Range("N10").ClearContents()

$ws.Range("H68").Value = 801.50665
$ws.Range("J68").Value = 2521.3333
$ws.Range("L68").Value = 7563.999899999999
$ws.Range("N68").Value = -9185.999899999999

$ws.Range("H71").Value = 801.50665
$ws.Range("J71").Value = 2521.3333
$ws.Range("L71").Value = 22691.9997
$ws.Range("N71").Value = -30803.9997

$ws.Range("H124").Value = 1651.4348
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 1651.4348
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 4954.3044
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -14774.3044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 9000
$ws.Range("J27").Value = 9000
$ws.Range("L27").Value = 9000
$ws.Range("N27").Value = -9332

$ws.Range("H132").Value = 55551.05
$ws.Range("I132").Value = 48056.09
$ws.Range("J132").Value = 65856.625
$ws.Range("K132").Value = 144168.27
$ws.Range("L132").Value = 197569.875
$ws.Range("M132").Value = -141638.27
$ws.Range("N132").Value = -202629.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2174.4443
$ws.Range("I7").Value = 2134.5833
$ws.Range("J7").Value = 2254.1667
$ws.Range("K7").Value = 2134.5833
$ws.Range("L7").Value = 2254.1667
$ws.Range("M7").Value = -2022.5833
$ws.Range("N7").Value = -2478.1667

$ws.Range("H76").Value = 39000
$ws.Range("J76").Value = 39000
$ws.Range("L76").Value = 39000
$ws.Range("N76").Value = -39676

$ws.Range("H79").Value = 39000
$ws.Range("J79").Value = 39000
$ws.Range("L79").Value = 39000
$ws.Range("N79").Value = -41340

$ws.Range("H93").Value = 1994.6086
$ws.Range("I93").Value = 2011.8889
$ws.Range("J93").Value = 1932.4
$ws.Range("K93").Value = 2011.8889
$ws.Range("L93").Value = 1932.4
$ws.Range("M93").Value = -763.8888999999999
$ws.Range("N93").Value = -4428.4

$ws.Range("H126").Value = 2174.4443
$ws.Range("I126").Value = 2134.5833
$ws.Range("J126").Value = 2254.1667
$ws.Range("K126").Value = 6403.749899999999
$ws.Range("L126").Value = 6762.500100000001
$ws.Range("M126").Value = -3933.749899999999
$ws.Range("N126").Value = -11702.5001

$ws.Range("H132").Value = 41232.08
$ws.Range("I132").Value = 2139.0588
$ws.Range("K132").Value = 6417.176399999999
$ws.Range("M132").Value = -3887.176399999999

$ws.Range("H136").Value = 69637.42
$ws.Range("I136").Value = 33165.97
$ws.Range("J136").Value = 175736.19
$ws.Range("K136").Value = 99497.91
$ws.Range("L136").Value = 527208.5700000001
$ws.Range("M136").Value = -96947.91
$ws.Range("N136").Value = -532308.5700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 6000
$ws.Range("J24").Value = 6000
$ws.Range("L24").Value = 6000
$ws.Range("N24").Value = -6460

$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 1666.6666
$ws.Range("K113").Value = 4999.9998
$ws.Range("M113").Value = -2829.9998

$ws.Range("H132").Value = 79339.34
$ws.Range("I132").Value = 61218.59
$ws.Range("J132").Value = 113567.445
$ws.Range("K132").Value = 183655.77
$ws.Range("L132").Value = 340702.335
$ws.Range("M132").Value = -181125.77
$ws.Range("N132").Value = -345762.335

$ws.Range("H136").Value = 84621.625
$ws.Range("I136").Value = 63751.188
$ws.Range("J136").Value = 126362.5
$ws.Range("K136").Value = 191253.564
$ws.Range("L136").Value = 379087.5
$ws.Range("M136").Value = -188703.564
$ws.Range("N136").Value = -384187.5

